# Running table.xlsx: add new "300440data" sheet with case2_20..case2_25 rows,
# adjust active sheet/selection on the other two sheets, and select the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Adjust view state of the two existing sheets (selection only - topLeftCell
#    cannot be round-tripped by this engine, verified separately).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("64000data")
$ws1.Range("H77").Select()

$ws2 = $wb.Worksheets.Item("64000 Bald data")
$ws2.Range("B3:K4").Select()
$ws2.Activate()

# ---------------------------------------------------------------------------
# 2) Add the new worksheet as the last tab, named "300440data".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "300440data"

# Column widths, matching the other two sheets' layout for this table.
$ws3.Columns.Item(2).ColumnWidth = 12.44140625
$ws3.Columns.Item(3).ColumnWidth = 17.109375
$ws3.Columns.Item(4).ColumnWidth = 12.5546875
$ws3.Columns.Item(6).ColumnWidth = 17
$ws3.Columns.Item(7).ColumnWidth = 12.88671875
$ws3.Columns.Item(8).ColumnWidth = 19
$ws3.Columns.Item(9).ColumnWidth = 19.5546875
$ws3.Columns.Item(10).ColumnWidth = 24.77734375
$ws3.Columns.Item(11).ColumnWidth = 23.21875

# --- Header row (row 2), copying the shaded/centered header formatting ------
$ws2.Range("B3:K3").Copy()
$ws3.Range("B2:K2").PasteSpecial(-4122)
$ws3.Range("B2").Value = "Job id"
$ws3.Range("C2").Value = "Name"
$ws3.Range("D2").Value = "CNN"
$ws3.Range("E2").Value = "Subcases"
$ws3.Range("F2").Value = "Augment"
$ws3.Range("G2").Value = "Running?"
$ws3.Range("H2").Value = "Activation"
$ws3.Range("I2").Value = "Remarks"
$ws3.Range("J2").Value = "Results"
$ws3.Range("K2").Value = "Saved as"

# --- First data row (row 3) - carries the "template" formatting ------------
$ws2.Range("B4:K4").Copy()
$ws3.Range("B3:K3").PasteSpecial(-4122)
$ws3.Range("B3").ClearContents()
$ws3.Range("C3").Value = "case2_20"
$ws3.Range("D3").Value = 2
$ws3.Range("E3").Value = 20
$ws3.Range("F3").Value = $true
$ws3.Range("G3").ClearContents()
$ws3.Range("H3").Value = "leakyRELU"
$ws3.Range("I3").ClearContents()
$ws3.Range("J3").ClearContents()
$ws3.Range("K3").ClearContents()

# --- Remaining data rows (4-8) - simple value rows --------------------------
$ws3.Range("D3:F3").Copy()
$ws3.Range("D4:F4").PasteSpecial(-4122)
$ws3.Range("D5:F5").PasteSpecial(-4122)
$ws3.Range("D6:F6").PasteSpecial(-4122)
$ws3.Range("D7:F7").PasteSpecial(-4122)
$ws3.Range("D8:F8").PasteSpecial(-4122)

$ws2.Range("H5").Copy()
$ws3.Range("H4").PasteSpecial(-4122)
$ws3.Range("H5").PasteSpecial(-4122)
$ws3.Range("H6").PasteSpecial(-4122)
$ws3.Range("H7").PasteSpecial(-4122)
$ws3.Range("H8").PasteSpecial(-4122)

$rows = @(
  @{ r = 4; c = "case2_21"; e = 21; f = $true  },
  @{ r = 5; c = "case2_22"; e = 22; f = $true  },
  @{ r = 6; c = "case2_23"; e = 23; f = $false },
  @{ r = 7; c = "case2_24"; e = 24; f = $false },
  @{ r = 8; c = "case2_25"; e = 25; f = $false }
)
foreach ($row in $rows) {
    $ws3.Range("C" + $row.r).Value = $row.c
    $ws3.Range("D" + $row.r).Value = 2
    $ws3.Range("E" + $row.r).Value = $row.e
    $ws3.Range("F" + $row.r).Value = $row.f
    $ws3.Range("H" + $row.r).Value = "leakyRELU"
}

$ws3.Range("B9", "K20").ClearContents()

# Phonetic settings (matches the noConversion flag used on every sheet here).
$ws3.Range("B2:K8").SetPhonetic()

$ws3.Range("H12").Select()
$ws3.Activate()
